$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL: https://hl7.fr/fhir/fr/medication/... -> https://hl7.fr/ig/fhir/medication/...
$ws.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-reconciliation-resolution"

# Name: FrMedicationReconciliationResolution -> FRMedicationReconciliationResolution
$ws.Range("B4").Value = "FRMedicationReconciliationResolution"

# Title: InterOp'Santé -> Interop'Santé
$ws.Range("B5").Value = "code system Interop'Santé - Résolution d'une divergence sur une ligne de traitement d'une FCT"

# Date: 2025-04-10T15:35:36+00:00 -> 2026-01-15T08:54:26+00:00
$ws.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction: "" -> FRANCE
$ws.Range("B11").Value = "FRANCE"
